# "Minor updates and tests"
#
# The underlying commit re-touches the centre/centre-vertical alignment
# already applied to the Wallcut table (columns C:D) - a cosmetic
# re-application that leaves the rendered formatting unchanged - and then
# leaves the sheet scrolled/selected at a different cell (P23) than where
# it had been (I12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# xlCenter = -4108
$xlCenter = -4108

# Re-apply the existing alignment on the data rows (C2:D39 -> centered
# horizontally and vertically) and on the merged header cell (C1:D1 ->
# centered horizontally only).
$ws.Range("C2:D39").HorizontalAlignment = $xlCenter
$ws.Range("C2:D39").VerticalAlignment = $xlCenter
$ws.Range("C1:D1").HorizontalAlignment = $xlCenter

# Move the active selection to P23, which is what was left selected when
# the workbook was saved.
$ws.Range("P23").Select()
